# Update the "2019 K-Means Clustering Data" table (slide 34) with the
# refreshed cluster-mean figures.
#
# Table layout (row 1 = header):
#   Cluster | Teams | MeanX2014HF | MeanX2015HF | MeanX2016HF | MeanX2017HF |
#   MeanX2018HF | MeanX2019HF | MeanPopulation | MeanSalary | MeanESPNRating
#
# Row 2 = Cluster 1, Row 3 = Cluster 2. Only the numeric values change;
# formatting/fills are left as-is.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(34)
$sh = $s.Shapes.Item(2)
$tbl = $sh.Table

# --- Cluster 1 (row 2) ---
$tbl.Cell(2,2).Shape.TextFrame.TextRange.Text  = "36"    # Teams: 30 -> 36
$tbl.Cell(2,3).Shape.TextFrame.TextRange.Text  = "1.9"   # MeanX2014HF: 1.5 -> 1.9
$tbl.Cell(2,4).Shape.TextFrame.TextRange.Text  = "2.3"   # MeanX2015HF: 1.9 -> 2.3
$tbl.Cell(2,5).Shape.TextFrame.TextRange.Text  = "2.3"   # MeanX2016HF: 2.4 -> 2.3
$tbl.Cell(2,6).Shape.TextFrame.TextRange.Text  = "1.8"   # MeanX2017HF: 2.5 -> 1.8
$tbl.Cell(2,7).Shape.TextFrame.TextRange.Text  = "1.5"   # MeanX2018HF: 2.1 -> 1.5
$tbl.Cell(2,8).Shape.TextFrame.TextRange.Text  = "1.1"   # MeanX2019HF: 1.5 -> 1.1
$tbl.Cell(2,9).Shape.TextFrame.TextRange.Text  = "1.8"   # MeanPopulation: 1.6 -> 1.8
$tbl.Cell(2,10).Shape.TextFrame.TextRange.Text = "14.0"  # MeanSalary: 12.7 -> 14.0
$tbl.Cell(2,11).Shape.TextFrame.TextRange.Text = "3.8"   # MeanESPNRating: 4.0 -> 3.8

# --- Cluster 2 (row 3) ---
$tbl.Cell(3,2).Shape.TextFrame.TextRange.Text  = "87"    # Teams: 93 -> 87
$tbl.Cell(3,3).Shape.TextFrame.TextRange.Text  = "0.5"   # MeanX2014HF: 0.7 -> 0.5
$tbl.Cell(3,4).Shape.TextFrame.TextRange.Text  = "0.3"   # MeanX2015HF: 0.6 -> 0.3
$tbl.Cell(3,6).Shape.TextFrame.TextRange.Text  = "0.6"   # MeanX2017HF: 0.4 -> 0.6
$tbl.Cell(3,7).Shape.TextFrame.TextRange.Text  = "0.7"   # MeanX2018HF: 0.6 -> 0.7
$tbl.Cell(3,8).Shape.TextFrame.TextRange.Text  = "0.9"   # MeanX2019HF: 0.7 -> 0.9
$tbl.Cell(3,10).Shape.TextFrame.TextRange.Text = "16.7"  # MeanSalary: 16.9 -> 16.7
